# Update "想去人数" (number of people interested) figures that changed
# between data pulls, on both the "展览" sheet and the aggregated
# "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows 2,3,4 -> column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 941
$wsExhibit.Range("F3").Value = 1800
$wsExhibit.Range("F4").Value = 406

# Sheet "全部类型" (All types) - same events appear in rows 4,5,6 -> column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 941
$wsAll.Range("F5").Value = 1800
$wsAll.Range("F6").Value = 406
